$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 485.61017
$ws.Range("J17").Value = 485.61017
$ws.Range("L17").Value = 1456.83051
$ws.Range("N17").Value = -1792.83051
# Row 62
$ws.Range("H62").Value = 7163.0454
$ws.Range("I62").Value = 9813.462
$ws.Range("J62").Value = 3334.6667
$ws.Range("K62").Value = 9813.462
$ws.Range("L62").Value = 3334.6667
$ws.Range("M62").Value = -9189.462
$ws.Range("N62").Value = -4582.6667
# Row 65
$ws.Range("H65").Value = 7163.0454
$ws.Range("I65").Value = 9813.462
$ws.Range("J65").Value = 3334.6667
$ws.Range("K65").Value = 49067.31
$ws.Range("L65").Value = 16673.3335
$ws.Range("M65").Value = -45947.31
$ws.Range("N65").Value = -22913.3335
# Row 74
$ws.Range("H74").Value = 4185.8945
$ws.Range("I74").Value = 4218
$ws.Range("J74").Value = 4141.75
$ws.Range("K74").Value = 4218
$ws.Range("L74").Value = 4141.75
$ws.Range("M74").Value = -3282
$ws.Range("N74").Value = -6013.75
# Row 77
$ws.Range("H77").Value = 4185.8945
$ws.Range("I77").Value = 4218
$ws.Range("J77").Value = 4141.75
$ws.Range("K77").Value = 21090
$ws.Range("L77").Value = 20708.75
$ws.Range("M77").Value = -16410
$ws.Range("N77").Value = -30068.75
# Row 116
$ws.Range("H116").Value = 2669.75
$ws.Range("I116").Value = 2420
$ws.Range("J116").Value = 2819.6
$ws.Range("K116").Value = 2420
$ws.Range("L116").Value = 2819.6
$ws.Range("M116").Value = 1022
$ws.Range("N116").Value = -9703.6
# Row 120
$ws.Range("H120").Value = 17833.666
$ws.Range("J120").Value = 17833.666
$ws.Range("L120").Value = 17833.666
$ws.Range("N120").Value = -27509.666
# Row 132
$ws.Range("H132").Value = 13896831
$ws.Range("I132").Value = 15630966
$ws.Range("K132").Value = 46892898
$ws.Range("M132").Value = -46890368

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 46848.363
$ws.Range("I45").Value = 200852
$ws.Range("J45").Value = 1553.1765
$ws.Range("K45").Value = 200852
$ws.Range("L45").Value = 1553.1765
$ws.Range("M45").Value = -200475
$ws.Range("N45").Value = -2307.1765
# Row 59
$ws.Range("H59").Value = 47780
$ws.Range("J59").Value = 47780
$ws.Range("L59").Value = 47780
$ws.Range("N59").Value = -49388
# Row 74
$ws.Range("H74").Value = 1264.7916
$ws.Range("I74").Value = 1449.4706
$ws.Range("J74").Value = 816.2857
$ws.Range("K74").Value = 1449.4706
$ws.Range("L74").Value = 816.2857
$ws.Range("M74").Value = -575.4706000000001
$ws.Range("N74").Value = -2564.2857
# Row 77
$ws.Range("H77").Value = 1264.7916
$ws.Range("I77").Value = 1449.4706
$ws.Range("J77").Value = 816.2857
$ws.Range("K77").Value = 7247.353000000001
$ws.Range("L77").Value = 4081.4285
$ws.Range("M77").Value = -2879.353000000001
$ws.Range("N77").Value = -12817.4285
# Row 97
$ws.Range("H97").Value = 2362.5
$ws.Range("I97").Value = 2464
$ws.Range("J97").Value = 2193.3333
$ws.Range("K97").Value = 2464
$ws.Range("L97").Value = 2193.3333
$ws.Range("M97").Value = -1968
$ws.Range("N97").Value = -3185.3333

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 2016.5883
$ws.Range("I99").Value = 1197
$ws.Range("J99").Value = 2590.3
$ws.Range("K99").Value = 1197
$ws.Range("L99").Value = 2590.3
$ws.Range("M99").Value = 301
$ws.Range("N99").Value = -5586.3

$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 3264.3333
$ws.Range("I62").Value = 2566.5
$ws.Range("J62").Value = 4660
$ws.Range("K62").Value = 2566.5
$ws.Range("L62").Value = 4660
$ws.Range("M62").Value = -1942.5
$ws.Range("N62").Value = -5908
# Row 65
$ws.Range("H65").Value = 3264.3333
$ws.Range("I65").Value = 2566.5
$ws.Range("J65").Value = 4660
$ws.Range("K65").Value = 12832.5
$ws.Range("L65").Value = 23300
$ws.Range("M65").Value = -9712.5
$ws.Range("N65").Value = -29540

$ws = $wb.Worksheets.Item("CUL")
# Row 11
$ws.Range("H11").Value = 50500
$ws.Range("I11").Value = 500
$ws.Range("J11").Value = 75500
$ws.Range("K11").Value = 1500
$ws.Range("L11").Value = 226500
$ws.Range("M11").Value = -1360
$ws.Range("N11").Value = -226780
# Row 46
$ws.Range("H46").Value = 300
$ws.Range("I46").Value = 300
$ws.Range("J46").Value = 300
$ws.Range("K46").Value = 900
$ws.Range("L46").Value = 900
$ws.Range("M46").Value = -809
$ws.Range("N46").Value = -1082
# Row 52
$ws.Range("H52").Value = 2606.6
$ws.Range("J52").Value = 2606.6
$ws.Range("L52").Value = 7819.799999999999
$ws.Range("N52").Value = -8351.799999999999
# Row 88
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
# Row 91
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
# Row 113
$ws.Range("H113").Value = 12536145
$ws.Range("I113").Value = 568.875
$ws.Range("J113").Value = 21256546
$ws.Range("K113").Value = 1706.625
$ws.Range("L113").Value = 63769638
$ws.Range("M113").Value = 463.375
$ws.Range("N113").Value = -63773978
# Row 131
$ws.Range("H131").Value = 766.08
$ws.Range("I131").Value = 425.75
$ws.Range("J131").Value = 812.48865
$ws.Range("K131").Value = 1277.25
$ws.Range("L131").Value = 2437.46595
$ws.Range("M131").Value = 3762.75
$ws.Range("N131").Value = -12517.46595

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 1985.5714
$ws.Range("I97").Value = 1985.5714
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1985.5714
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1489.5714
$ws.Range("N97").ClearContents()
# Row 102
$ws.Range("H102").Value = 1384.0834
$ws.Range("I102").Value = 1305.9
$ws.Range("J102").Value = 1775
$ws.Range("K102").Value = 1305.9
$ws.Range("L102").Value = 1775
$ws.Range("M102").Value = 316.0999999999999
$ws.Range("N102").Value = -5019
# Row 122
$ws.Range("H122").Value = 31257696
$ws.Range("I122").Value = 38470324
$ws.Range("J122").Value = 2969.3333
$ws.Range("K122").Value = 115410972
$ws.Range("L122").Value = 8907.999899999999
$ws.Range("M122").Value = -115408522
$ws.Range("N122").Value = -13807.9999

$ws = $wb.Worksheets.Item("LTW")
# Row 60
$ws.Range("H60").Value = 17000
$ws.Range("J60").Value = 17000
$ws.Range("L60").Value = 17000
$ws.Range("N60").Value = -18018
# Row 93
$ws.Range("H93").Value = 1128.3636
$ws.Range("I93").Value = 1282.4
$ws.Range("K93").Value = 1282.4
$ws.Range("M93").Value = -34.40000000000009
# Row 100
$ws.Range("H100").Value = 2593.5151
$ws.Range("I100").Value = 3700.75
$ws.Range("J100").Value = 2440.7932
$ws.Range("K100").Value = 3700.75
$ws.Range("L100").Value = 2440.7932
$ws.Range("M100").Value = -3159.75
$ws.Range("N100").Value = -3522.7932
# Row 132
$ws.Range("H132").Value = 6796.8374
$ws.Range("I132").Value = 1568.9259
$ws.Range("J132").Value = 15618.9375
$ws.Range("K132").Value = 4706.7777
$ws.Range("L132").Value = 46856.8125
$ws.Range("M132").Value = -2176.7777
$ws.Range("N132").Value = -51916.8125
# Row 136
$ws.Range("H136").Value = 39411324
$ws.Range("I136").Value = 6496501
$ws.Range("J136").Value = 142857900
$ws.Range("K136").Value = 19489503
$ws.Range("L136").Value = 428573700
$ws.Range("M136").Value = -19486953
$ws.Range("N136").Value = -428578800

$ws = $wb.Worksheets.Item("WVR")
# Row 64
$ws.Range("H64").Value = 10944
$ws.Range("J64").Value = 10944
$ws.Range("L64").Value = 10944
$ws.Range("N64").Value = -11440
# Row 67
$ws.Range("H67").Value = 10944
$ws.Range("J67").Value = 10944
$ws.Range("L67").Value = 10944
$ws.Range("N67").Value = -12660
# Row 136
$ws.Range("H136").Value = 2588.375
$ws.Range("I136").Value = 1667
$ws.Range("K136").Value = 5001
$ws.Range("M136").Value = -2451
